$d = $word.ActiveDocument

# Table 1 (GEOMETRIJA overview table), "Prostorne relacije" row:
# "1, 2, 6, 7, 8, 14, 15, 45, 77, 83" -> "1, 2, 6, 7, 8, 14, 15, 77, 83"
# (drop "45, " from the list of task numbers)
$d.Content.Find.Execute("1, 2, 6, 7, 8, 14, 15, 45, 77", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "1, 2, 6, 7, 8, 14, 15, 77", 2) | Out-Null

# Same table, "Spoljašnja, unutrašnja oblast" row:
# "9, 56" -> "87, 56"
$d.Content.Find.Execute("9, 56", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "87, 56", 2) | Out-Null

# Same table, last row ("Merenje dužine i težine"):
# "45, 46, 47" -> "45, 46, 47, 108"
$d.Content.Find.Execute("45, 46, 47", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "45, 46, 47, 108", 2) | Out-Null
